$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the weekly data block (row 182),
# pushing the existing rows (182..232) down to (184..234).
$ws.Rows("182:183").Insert()

# New week's data (first new row) - "Primera" quality
$ws.Range("A182").Value = 11
$ws.Range("B182").Value = "Vega Monumental Concepción"
$ws.Range("C182").Value = "Bíobío"
$ws.Range("D182").Value = 44736
$ws.Range("E182").Value = 8
$ws.Range("F182").Value = 100114013
$ws.Range("G182").Value = "Zanahoria"
$ws.Range("H182").Value = "Sin especificar"
$ws.Range("I182").Value = "Primera"
$ws.Range("J182").Value = 1000
$ws.Range("K182").Value = 5500
$ws.Range("L182").Value = 6000
$ws.Range("M182").Value = 5750
$ws.Range("N182").Value = "$/saco 20 kilos"
$ws.Range("O182").Value = "Región de Ñuble"
$ws.Range("P182").Value = 288
$ws.Range("Q182").Value = 20
$ws.Range("R182").Value = "Hortaliza"

# New week's data (second new row) - "Segunda" quality
$ws.Range("A183").Value = 11
$ws.Range("B183").Value = "Vega Monumental Concepción"
$ws.Range("C183").Value = "Bíobío"
$ws.Range("D183").Value = 44736
$ws.Range("E183").Value = 8
$ws.Range("F183").Value = 100114013
$ws.Range("G183").Value = "Zanahoria"
$ws.Range("H183").Value = "Sin especificar"
$ws.Range("I183").Value = "Segunda"
$ws.Range("J183").Value = 500
$ws.Range("K183").Value = 5000
$ws.Range("L183").Value = 5000
$ws.Range("M183").Value = 5000
$ws.Range("N183").Value = "$/saco 20 kilos"
$ws.Range("O183").Value = "Región de Ñuble"
$ws.Range("P183").Value = 250
$ws.Range("Q183").Value = 20
$ws.Range("R183").Value = "Hortaliza"
